# Updates the cryptos price list (columns D = Price, E = Volume(1h))
# for rows 2..51 to reflect the latest scrape, per commit
# "Updated cryptos list on Fri Aug  2 21:51:58 UTC 2024 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; DNew = "62.010.36"; DForceText = $false; ENew = "  -4.70%  " },
    @{ Row = 3; DNew = "2.999.83"; DForceText = $false; ENew = "  -5.71%  " },
    @{ Row = 4; DNew = "1.00"; DForceText = $true; ENew = "  -0.01%  " },
    @{ Row = 5; DNew = "550.39"; DForceText = $true; ENew = "  -4.03%  " },
    @{ Row = 6; DNew = "154.32"; DForceText = $true; ENew = "  -7.71%  " },
    @{ Row = 7; DNew = "1.00"; DForceText = $true; ENew = "  -0.07%  " },
    @{ Row = 8; DNew = "0.566"; DForceText = $true; ENew = "  -4.92%  " },
    @{ Row = 9; DNew = "3.008.44"; DForceText = $false; ENew = "  -5.34%  " },
    @{ Row = 10; DNew = $null; DForceText = $false; ENew = "  -5.61%  " },
    @{ Row = 11; DNew = "6.26"; DForceText = $true; ENew = "  -6.57%  " },
    @{ Row = 12; DNew = $null; DForceText = $false; ENew = "  -5.33%  " },
    @{ Row = 13; DNew = "3.523.38"; DForceText = $false; ENew = "  -5.56%  " },
    @{ Row = 14; DNew = $null; DForceText = $false; ENew = "  -3.86%  " },
    @{ Row = 15; DNew = "62.069.29"; DForceText = $false; ENew = "  -5.18%  " },
    @{ Row = 16; DNew = "23.75"; DForceText = $true; ENew = "  -7.38%  " },
    @{ Row = 17; DNew = "3.000.86"; DForceText = $false; ENew = "  -5.40%  " },
    @{ Row = 18; DNew = $null; DForceText = $false; ENew = "  -5.52%  " },
    @{ Row = 19; DNew = "393.50"; DForceText = $true; ENew = "  -5.60%  " },
    @{ Row = 20; DNew = $null; DForceText = $false; ENew = "  -3.35%  " },
    @{ Row = 21; DNew = "11.98"; DForceText = $true; ENew = "  -6.06%  " },
    @{ Row = 22; DNew = $null; DForceText = $false; ENew = "  -7.20%  " },
    @{ Row = 23; DNew = $null; DForceText = $false; ENew = "  -0.26%  " },
    @{ Row = 24; DNew = "65.14"; DForceText = $true; ENew = "  -5.14%  " },
    @{ Row = 25; DNew = $null; DForceText = $false; ENew = "  -4.21%  " },
    @{ Row = 26; DNew = "0.186"; DForceText = $true; ENew = "  -7.95%  " },
    @{ Row = 27; DNew = $null; DForceText = $false; ENew = "  -9.60%  " },
    @{ Row = 28; DNew = "0.999"; DForceText = $true; ENew = "  -0.16%  " },
    @{ Row = 29; DNew = "8.51"; DForceText = $true; ENew = "  -4.36%  " },
    @{ Row = 30; DNew = $null; DForceText = $false; ENew = "  -0.02%  " },
    @{ Row = 31; DNew = $null; DForceText = $false; ENew = "  -5.70%  " },
    @{ Row = 32; DNew = $null; DForceText = $false; ENew = "  -4.66%  " },
    @{ Row = 33; DNew = "159.55"; DForceText = $true; ENew = "  +2.29%  " },
    @{ Row = 34; DNew = $null; DForceText = $false; ENew = "  -6.71%  " },
    @{ Row = 35; DNew = $null; DForceText = $false; ENew = "  -5.37%  " },
    @{ Row = 36; DNew = $null; DForceText = $false; ENew = "  -5.10%  " },
    @{ Row = 37; DNew = "1.29"; DForceText = $true; ENew = "  -5.75%  " },
    @{ Row = 38; DNew = $null; DForceText = $false; ENew = "  -9.47%  " },
    @{ Row = 39; DNew = "2.452.61"; DForceText = $false; ENew = "  -10.01%  " },
    @{ Row = 40; DNew = "3.93"; DForceText = $true; ENew = "  -4.70%  " },
    @{ Row = 41; DNew = "22.46"; DForceText = $true; ENew = "  -6.03%  " },
    @{ Row = 42; DNew = "37.20"; DForceText = $true; ENew = "  -4.62%  " },
    @{ Row = 43; DNew = "0.664"; DForceText = $true; ENew = "  -6.40%  " },
    @{ Row = 44; DNew = $null; DForceText = $false; ENew = "  -6.38%  " },
    @{ Row = 45; DNew = "0.999"; DForceText = $true; ENew = "  -0.15%  " },
    @{ Row = 46; DNew = $null; DForceText = $false; ENew = "  -5.97%  " },
    @{ Row = 47; DNew = "4.96"; DForceText = $true; ENew = "  -11.27%  " },
    @{ Row = 48; DNew = $null; DForceText = $false; ENew = "  -7.29%  " },
    @{ Row = 49; DNew = "0.0949"; DForceText = $true; ENew = "  -3.96%  " },
    @{ Row = 50; DNew = $null; DForceText = $false; ENew = "  +0.41%  " },
    @{ Row = 51; DNew = "266.43"; DForceText = $true; ENew = "  -9.28%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.DNew) {
        $cell = $ws.Cells.Item($row, 4)
        if ($u.DForceText) {
            # The new price string (e.g. "1.00", "37.20") is valid numeric
            # literal syntax, so Excel would silently coerce it to a number.
            # Force a text number format while assigning so it is stored as
            # a string, matching the source data, then restore the cell's
            # normal style so no stray formatting is left behind.
            $cell.NumberFormat = "@"
            $cell.Value = $u.DNew
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.DNew
        }
    }

    $ws.Cells.Item($row, 5).Value = $u.ENew
}
